{"js": "// HTH: Fixing a little detail to chapter 23\n// \"The studio was even well-ordered than the atrium ...\"\n//   -> \"The studio was well-ordered like the atrium ...\"\n// i.e. drop \"even \" and swap \"than\" -> \"like\".\n\nconst searchResults = context.document.body.search(\n  \"The studio was even well-ordered than the atrium\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target sentence fragment not found\");\n}\n\nconst target = searchResults.items[0];\ntarget.insertText(\n  \"The studio was well-ordered like the atrium\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# HTH: Fixing a little detail to chapter 23\n# \"The studio was even well-ordered than the atrium ...\"\n#   -> \"The studio was well-ordered like the atrium ...\"\n# i.e. drop \"even \" and swap \"than\" -> \"like\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"The studio was even well-ordered than the atrium\",  # FindText\n    $true,                                                # MatchCase\n    $false,                                               # MatchWholeWord\n    $false,                                               # MatchWildcards\n    $false,                                               # MatchSoundsLike\n    $false,                                               # MatchAllWordForms\n    $true,                                                # Forward\n    1,                                                     # Wrap (wdFindContinue)\n    $false,                                               # Format\n    \"The studio was well-ordered like the atrium\",       # ReplaceWith\n    2                                                      # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw \"Target sentence fragment not found\"\n}\n"}
